$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 663953.9399999999
$ws.Range("J17").Value = 1427289.9
$ws.Range("L17").Value = 4281869.699999999
$ws.Range("N17").Value = -4282205.699999999
$ws.Range("H32").Value = 3055.1177
$ws.Range("I32").Value = 2378.75
$ws.Range("K32").Value = 2378.75
$ws.Range("M32").Value = -2052.75
$ws.Range("H51").Value = 2818.45
$ws.Range("I51").Value = 2473
$ws.Range("K51").Value = 2473
$ws.Range("M51").Value = -1989
$ws.Range("H58").Value = 10113.2
$ws.Range("I58").Value = 126
$ws.Range("K58").Value = 378
$ws.Range("M58").Value = -228
$ws.Range("H113").Value = 12788.777
$ws.Range("I113").Value = 7219.8
$ws.Range("K113").Value = 7219.8
$ws.Range("M113").Value = -3965.8
$ws.Range("H137").Value = 1999.4546
$ws.Range("I137").Value = 2099.6
$ws.Range("J137").Value = 1916
$ws.Range("K137").Value = 6298.799999999999
$ws.Range("L137").Value = 5748
$ws.Range("M137").Value = -3748.799999999999
$ws.Range("N137").Value = -10848
$ws.Range("H138").Value = 7695075
$ws.Range("I138").Value = 1256.5454
$ws.Range("K138").Value = 3769.6362
$ws.Range("M138").Value = 1370.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 87.3
$ws.Range("I5").Value = 74.59999999999999
$ws.Range("K5").Value = 74.59999999999999
$ws.Range("M5").Value = 37.40000000000001
$ws.Range("H61").Value = 3757.1606
$ws.Range("I61").Value = 2590.0286
$ws.Range("K61").Value = 2590.0286
$ws.Range("M61").Value = -2378.0286
$ws.Range("H74").Value = 5059.5
$ws.Range("I74").Value = 1217.4595
$ws.Range("K74").Value = 1217.4595
$ws.Range("M74").Value = -343.4594999999999
$ws.Range("H77").Value = 5059.5
$ws.Range("I77").Value = 1217.4595
$ws.Range("K77").Value = 6087.2975
$ws.Range("M77").Value = -1719.2975
$ws.Range("H88").Value = 2814.4285
$ws.Range("I88").Value = 2006
$ws.Range("K88").Value = 2006
$ws.Range("M88").Value = -1600
$ws.Range("H91").Value = 2814.4285
$ws.Range("I91").Value = 2006
$ws.Range("K91").Value = 2006
$ws.Range("M91").Value = -602
$ws.Range("H122").Value = 2596.2778
$ws.Range("I122").Value = 2703.3
$ws.Range("J122").Value = 2462.5
$ws.Range("K122").Value = 8109.900000000001
$ws.Range("L122").Value = 7387.5
$ws.Range("M122").Value = -5659.900000000001
$ws.Range("N122").Value = -12287.5
$ws.Range("H133").Value = 96725.7
$ws.Range("J133").Value = 96725.7
$ws.Range("L133").Value = 96725.7
$ws.Range("N133").Value = -101785.7
$ws.Range("H135").Value = 33999.75
$ws.Range("J135").Value = 33999.75
$ws.Range("L135").Value = 33999.75
$ws.Range("N135").Value = -44139.75
$ws.Range("H136").Value = 3757.1606
$ws.Range("I136").Value = 2590.0286
$ws.Range("K136").Value = 7770.085800000001
$ws.Range("M136").Value = -5220.085800000001
$ws.Range("H139").Value = 92949.39999999999
$ws.Range("J139").Value = 92949.39999999999
$ws.Range("L139").Value = 92949.39999999999
$ws.Range("N139").Value = -103229.4
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 87.3
$ws.Range("I4").Value = 74.59999999999999
$ws.Range("K4").Value = 74.59999999999999
$ws.Range("M4").Value = 40.40000000000001
$ws.Range("H132").Value = 77166.336
$ws.Range("J132").Value = 77166.336
$ws.Range("L132").Value = 77166.336
$ws.Range("N132").Value = -87286.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 42170
$ws.Range("I31").Value = 49640.523
$ws.Range("J31").Value = 2949.75
$ws.Range("K31").Value = 49640.523
$ws.Range("L31").Value = 2949.75
$ws.Range("M31").Value = -49345.523
$ws.Range("N31").Value = -3539.75
$ws.Range("H34").Value = 42170
$ws.Range("I34").Value = 49640.523
$ws.Range("J34").Value = 2949.75
$ws.Range("K34").Value = 49640.523
$ws.Range("L34").Value = 2949.75
$ws.Range("M34").Value = -49438.523
$ws.Range("N34").Value = -3353.75
$ws.Range("H41").Value = 40499.25
$ws.Range("I41").Value = 5999.5
$ws.Range("K41").Value = 5999.5
$ws.Range("M41").Value = -5571.5
$ws.Range("H103").Value = 9818.583000000001
$ws.Range("I103").Value = 7983.909
$ws.Range("K103").Value = 7983.909
$ws.Range("M103").Value = -6811.909
$ws.Range("H122").Value = 1854.4286
$ws.Range("I122").Value = 1596.4
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 4789.200000000001
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -2339.200000000001
$ws.Range("N122").Value = -12398.5
$ws.Range("H132").Value = 4646.8
$ws.Range("I132").Value = 4646.8
$ws.Range("K132").Value = 13940.4
$ws.Range("M132").Value = -11410.4
$ws.Range("H134").Value = 19609.045
$ws.Range("I134").Value = 8069.95
$ws.Range("J134").Value = 135000
$ws.Range("K134").Value = 24209.85
$ws.Range("L134").Value = 405000
$ws.Range("M134").Value = -21674.85
$ws.Range("N134").Value = -410070
$ws.Range("H99").Value = 4875
$ws.Range("I99").Value = 4500
$ws.Range("J99").Value = 5250
$ws.Range("K99").Value = 4500
$ws.Range("L99").Value = 5250
$ws.Range("M99").Value = -3002
$ws.Range("N99").Value = -8246
$ws.Range("H126").Value = 4875
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 5250
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 15750
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -20690

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 74.5
$ws.Range("I38").Value = 63.42857
$ws.Range("J38").Value = 85.57143000000001
$ws.Range("K38").Value = 190.28571
$ws.Range("L38").Value = 256.71429
$ws.Range("M38").Value = 156.71429
$ws.Range("N38").Value = -950.71429
$ws.Range("H58").Value = 999.5
$ws.Range("I58").Value = 999
$ws.Range("K58").Value = 2997
$ws.Range("M58").Value = -2869
$ws.Range("H59").Value = 55
$ws.Range("I59").Value = 55
$ws.Range("K59").Value = 165
$ws.Range("M59").Value = 375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65.88
$ws.Range("I2").Value = 65.052635
$ws.Range("J2").Value = 68.5
$ws.Range("K2").Value = 65.052635
$ws.Range("L2").Value = 68.5
$ws.Range("M2").Value = 47.947365
$ws.Range("N2").Value = -294.5
$ws.Range("H24").Value = 15888.333
$ws.Range("H102").Value = 50002290
$ws.Range("I102").Value = 2409.8948
$ws.Range("K102").Value = 2409.8948
$ws.Range("M102").Value = -787.8948
$ws.Range("H107").Value = 261
$ws.Range("I107").Value = 266.14285
$ws.Range("K107").Value = 266.14285
$ws.Range("M107").Value = 1653.85715
$ws.Range("H122").Value = 2890.3125
$ws.Range("I122").Value = 2603.2144
$ws.Range("K122").Value = 7809.6432
$ws.Range("M122").Value = -5359.6432
$ws.Range("H69").Value = 15000
$ws.Range("J69").Value = 15000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16498
$ws.Range("H72").Value = 15000
$ws.Range("J72").Value = 15000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -52488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4819.625
$ws.Range("I7").Value = 4182.8423
$ws.Range("K7").Value = 4182.8423
$ws.Range("M7").Value = -4070.8423
$ws.Range("H40").Value = 6456.6
$ws.Range("J40").Value = 7624.6665
$ws.Range("L40").Value = 7624.6665
$ws.Range("N40").Value = -7896.6665
$ws.Range("H93").Value = 1421.5
$ws.Range("I93").Value = 1415
$ws.Range("J93").Value = 1441
$ws.Range("K93").Value = 1415
$ws.Range("L93").Value = 1441
$ws.Range("M93").Value = -167
$ws.Range("N93").Value = -3937
$ws.Range("H126").Value = 4819.625
$ws.Range("I126").Value = 4182.8423
$ws.Range("K126").Value = 12548.5269
$ws.Range("M126").Value = -10078.5269
